# Generate Report for Handback
# -----------------------------------------------------------------------
# This script reproduces, through Excel COM automation, the "handback"
# report-generation pass: the localization status moves from
# "Ready for handoff" to "Handed back: in sync with en-US", the target
# (translated) files + their handback timestamps are now known, and the
# relevant columns are widened so the longer values are readable.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# ColumnWidth is expressed in characters and Excel snaps it to a whole
# number of pixels (MaximumDigitWidth = 6 for the default Calibri 11
# used by this workbook): stored_width = Round(chars*6)/6 + 5/6.
# The two helper values below are the character widths that round-trip
# to the target stored widths of ~29.98 and 40 characters respectively.
$wMedium = 29.166666666666668   # -> stored width ~29.98 (status columns)
$wWide   = 39.166666666666664   # -> stored width 40     (file/name columns)

# ---------------------------------------------------------------------
# 1. Status: "Ready for handoff" -> "Handed back: in sync with en-US"
# ---------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# 2. Overview sheet: widen the per-locale status columns (E, F)
# ---------------------------------------------------------------------
$overview.Columns.Item(5).ColumnWidth = $wMedium
$overview.Columns.Item(6).ColumnWidth = $wMedium

# ---------------------------------------------------------------------
# 3. zh-cn / de-de sheets: record the handback target file + datetime,
#    widen Status (C), Latest Target File (H) and Latest Handback File
#    (I) to fit the now-populated values.
# ---------------------------------------------------------------------
$targetFileUrl = "https://github.com/OpenLocalizationTestOrg/oltest/blob/b6985c375a33e082d1943053e3553409c92c6197/e2e/dfd87fc3-e78e-4fa1-95a1-f875ed042ef9.md"
$targetFileName = "dfd87fc3-e78e-4fa1-95a1-f875ed042ef9.md"

foreach ($ws in @($zhcn, $dede)) {
    $ws.Columns.Item(3).ColumnWidth = $wMedium
    $ws.Columns.Item(8).ColumnWidth = $wWide
    $ws.Columns.Item(9).ColumnWidth = $wWide

    $ws.Hyperlinks.Add($ws.Range("H2"), $targetFileUrl, "", "", $targetFileName)
    $ws.Hyperlinks.Add($ws.Range("H3"), $targetFileUrl, "", "", $targetFileName)
}

# Latest Handback File: the generated xliff that was handed back.
$zhcn.Range("I2").Value = "dfd87fc3-e78e-4fa1-95a1-f875ed042ef9.2f07c35368b9579b291927cc6804ae3b6a7af3f0.zh-cn.xlf"
$zhcn.Range("I3").Value = "dfd87fc3-e78e-4fa1-95a1-f875ed042ef9.2f07c35368b9579b291927cc6804ae3b6a7af3f0.zh-cn.xlf"

$dede.Range("I2").Value = "dfd87fc3-e78e-4fa1-95a1-f875ed042ef9.2f07c35368b9579b291927cc6804ae3b6a7af3f0.de-de.xlf"
$dede.Range("I3").Value = "dfd87fc3-e78e-4fa1-95a1-f875ed042ef9.2f07c35368b9579b291927cc6804ae3b6a7af3f0.de-de.xlf"

# Latest Handback DateTime: zh-cn finished first, de-de a little later.
$zhcn.Range("J2").Value = "2016-07-26 08:10:36"
$zhcn.Range("J3").Value = "2016-07-26 08:10:36"

$dede.Range("J2").Value = "2016-07-26 08:10:51"
$dede.Range("J3").Value = "2016-07-26 08:10:51"
